# Commit: "Fruta / hortaliza, semanal"
# A new weekly price-record row is inserted into the daily price log for
# "Terminal Hortofrutícola Agro Chillán - Naranja" at sheet row 415 (just
# below the existing row 414), pushing every row that used to be 415-520
# down to 416-521. Dimension grows from A1:T520 to A1:T521.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row before row 415 -- everything currently at 415..520
# shifts down to 416..521 (Excel also expands the sheet dimension for us).
$ws.Rows("415:415").Insert()

# The columns that are constant for every data row in this sheet
# (Mercado ID / Mercado / Región / Codreg / Tipo / Producto ID / Producto /
# Categoría ID / Categoría) need to be filled in for the freshly inserted
# row too, since Insert() only shifted the rows below -- it left row 415
# blank.
$ws.Range("A415").Value = 7
$ws.Range("B415").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C415").Value = "Ñuble"
$ws.Range("D415").Value = 44943
$ws.Range("E415").Value = 16
$ws.Range("F415").Value = "Fruta"
$ws.Range("G415").Value = 100102
$ws.Range("H415").Value = "Cítricos"
$ws.Range("I415").Value = 100102005
$ws.Range("J415").Value = "Naranja"
$ws.Range("K415").Value = "Valencia"
$ws.Range("L415").Value = "Primera"
$ws.Range("M415").Value = 100
$ws.Range("N415").Value = 10000
$ws.Range("O415").Value = 11000
$ws.Range("P415").Value = 10500
$ws.Range("Q415").Value = "$/caja 15 kilos granel"
$ws.Range("R415").Value = "Región de O'Higgins"
$ws.Range("S415").Value = 700
$ws.Range("T415").Value = 15
